# Update TC_ID Excel SCD0017 until SCD0025 and Update TC_ID Solution SCD0006 until SCD0025
#
# This test-case sheet was renumbered from the old Jira id (DGS-316) to the
# new SCD naming scheme (SCD0018-024), and the worksheet/tab itself was
# renamed from "SCD0301" to "SCD0018" to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The tab name follows the TC_ID family, so rename it to match.
$ws.Name = "SCD0018"

# TC_ID column (B) - every data row (2-6) gets the new id.
$ws.Range("B2:B6").Value = "SCD0018-024"

# Column B needs to be a bit wider to fit the longer "SCD0018-024" text.
$ws.Columns.Item(2).ColumnWidth = 11.65

# Restore the scroll position / active cell from the authoring session
# (previously topLeftCell B4 / selection D6 -> now topLeftCell A4 / selection B7).
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B7").Select() | Out-Null
